# AddShift.xlsx refactor:
#  - Renumber the "Test Case" column (B) to be sequential (skips removed)
#  - Mark several rows' "Revise" column (A) from "Y" to "N"
#  - Fill in the "Revise result" column (J) to mirror/override the Result (I) column
#  - Update the sheet view (selection + zoom) and formulas recalc automatically

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A ("Revise"): switch selected rows from Y to N ---
$reviseNoRows = 3,4,7,8,10,12,13,14,16,17,18,19,20,21
foreach ($r in $reviseNoRows) {
    $ws.Range("A$r").Value = "N"
}

# --- Column B ("Test Case"): renumber sequentially ---
$bValues = @{
    2  = 1
    3  = 2
    4  = 3
    5  = 4
    6  = 5
    7  = 6
    8  = 7
    9  = 8
    10 = 9
    11 = 10
    12 = 11
    13 = 12
    14 = 13
    15 = 14
    16 = 15
    17 = 16
    18 = 17
    19 = 18
    20 = 19
    21 = 20
    22 = 21
}
foreach ($r in $bValues.Keys) {
    $ws.Range("B$r").Value = $bValues[$r]
}

# --- Column J ("Revise result"): fill in with Pass/Fail, mirroring column I
#     except rows 9 and 13 which were manually corrected to Pass ---
$jValues = @{
    2  = "Pass"
    3  = "Fail"
    4  = "Fail"
    5  = "Fail"
    6  = "Pass"
    7  = "Fail"
    8  = "Fail"
    9  = "Pass"
    10 = "Pass"
    11 = "Fail"
    12 = "Fail"
    13 = "Pass"
    14 = "Pass"
    15 = "Fail"
    16 = "Fail"
    17 = "Pass"
    18 = "Pass"
    19 = "Pass"
    20 = "Pass"
    21 = "Fail"
    22 = "Pass"
}
foreach ($r in $jValues.Keys) {
    $ws.Range("J$r").Value = $jValues[$r]
}

# --- Sheet view: clear the frozen/top-left cell scroll position, select K16, zoom to 86% ---
$ws.Activate()
$ws.Range("K16").Select()
$excel.ActiveWindow.Zoom = 86
try {
    # Best-effort: move the (off-screen) app window closer to the origin,
    # matching the saved workbook window position.
    $excel.ActiveWindow.Left = -108
} catch {
    # Not fatal if the host doesn't support positioning the window.
}

# --- Recalculate so the COUNTIF/TEXT summary cells (N21:O22) reflect column J ---
$excel.Calculate()
